$d = $word.ActiveDocument

# --- Change 1: Semaine 1 block ---
# Original paragraphs (in order):
#   "- Première chose à faire : prendre un temps pour s'approprier les datasets et voir ce qu'il est possible de faire. "
#   "- Définir l'objectif de l'entreprise, que propose t-elle comme service ?" (with spell-check proofErr runs)
#   "- Fixer le rôle de chacun dans l'entreprise"
# Target paragraphs (in order):
#   "- Définir l'objectif de l'entreprise, que propose t" + "'" + "elle comme service ?"
#   "- Fixer le rôle de chacun dans l'entreprise"
#   "- Prendre un temps pour s'approprier les datasets " + "et les notebooks déjà existants sur le sujet (disponible sur le Git)"

$anchor1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("- Première chose")) {
        $anchor1 = $p
        break
    }
}

$start1 = $anchor1.Range.Start
$p2 = $anchor1.Next()
$p3 = $p2.Next()
$end1 = $p3.Range.End

$r1 = $d.Range($start1, $end1)

$xml1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>- Définir l'objectif de l'entreprise, que propose t</w:t></w:r><w:r><w:t>’</w:t></w:r><w:r><w:t>elle comme service ?</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>- Fixer le rôle de chacun dans l'entreprise</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">- Prendre un temps pour s'approprier les datasets </w:t></w:r><w:r><w:t>et les notebooks déjà existants sur le sujet (disponible sur le Git)</w:t></w:r></w:p>
'@

$r1.InsertXML($xml1)

# --- Change 2: Semaine 2 block ---
# Original paragraph: "- Extraction des données pour faire une analyse descriptive des données sur un " + "Jupyter" (spell-check) + " Notebook"
# Target paragraph: single run "- Extraction des données pour faire une analyse descriptive des données sur un Jupyter Notebook"

$anchor2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("- Extraction des donn")) {
        $anchor2 = $p
        break
    }
}

$r2 = $anchor2.Range

$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>- Extraction des données pour faire une analyse descriptive des données sur un Jupyter Notebook</w:t></w:r></w:p>'

$r2.InsertXML($xml2)

Write-Output "done"
